$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump version + refresh date ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "1.8.2"
$meta.Range("B8").Value = "2023-09-01T14:45:29-04:00"

# --- Elements sheet: root "Extension" row now also carries the full
#     invariant text (ele-1 / ext-1) in its Invariants column (AJ),
#     same text already shown on the Extension.extension row (AJ3). ---
$elements = $wb.Worksheets.Item("Elements")
$invariantText = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`next-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"
$elements.Range("AJ1").Value = $invariantText
